$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the numeric-looking columns (D, E, G) as Text so Excel
# keeps the new values as strings (matching the source data) instead
# of auto-converting them to numbers/percentages on input.
$colD = $ws.Range("D2:D51")
$colE = $ws.Range("E2:E51")
$colG = $ws.Range("G2:G51")
$colD.NumberFormat = "@"
$colE.NumberFormat = "@"
$colG.NumberFormat = "@"

$ws.Range("D2").Value = "297.23"
$ws.Range("E2").Value = "2.12%"
$ws.Range("G2").Value = "4"
$ws.Range("D3").Value = "41.93"
$ws.Range("E3").Value = "3.46%"
$ws.Range("G3").Value = "4"
$ws.Range("D4").Value = "5.007"
$ws.Range("E4").Value = "-0.12%"
$ws.Range("G4").Value = "4"
$ws.Range("D5").Value = "0.07524"
$ws.Range("E5").Value = "2.64%"
$ws.Range("G5").Value = "4"
$ws.Range("B6").Value = "FTXToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D6").Value = "1.584"
$ws.Range("E6").Value = "3.55%"
$ws.Range("G6").Value = "4"
$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D7").Value = "0.9279"
$ws.Range("E7").Value = "0.59%"
$ws.Range("G7").Value = "4"
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").Value = "2.401"
$ws.Range("E8").Value = "1.61%"
$ws.Range("G8").Value = "4"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "0.1194"
$ws.Range("E9").Value = "-1.53%"
$ws.Range("G9").Value = "4"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "0.1831"
$ws.Range("E10").Value = "5.02%"
$ws.Range("G10").Value = "4"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "0.08920"
$ws.Range("E11").Value = "2.96%"
$ws.Range("G11").Value = "4"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "0.04077"
$ws.Range("E12").Value = "-4.90%"
$ws.Range("G12").Value = "4"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "0.1048"
$ws.Range("E13").Value = "-0.52%"
$ws.Range("G13").Value = "4"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "0.001288"
$ws.Range("E14").Value = "1.11%"
$ws.Range("G14").Value = "4"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "0.005977"
$ws.Range("E15").Value = "2.23%"
$ws.Range("G15").Value = "4"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "3.355"
$ws.Range("E16").Value = "0.57%"
$ws.Range("G16").Value = "4"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").Value = "4.380"
$ws.Range("E17").Value = "1.92%"
$ws.Range("G17").Value = "4"
$ws.Range("D18").Value = "0.3313"
$ws.Range("E18").Value = "0.77%"
$ws.Range("G18").Value = "4"
$ws.Range("D19").Value = "8.082"
$ws.Range("E19").Value = "4.67%"
$ws.Range("G19").Value = "4"
$ws.Range("D20").Value = "0.1410"
$ws.Range("E20").Value = "1.41%"
$ws.Range("G20").Value = "4"
$ws.Range("E21").Value = "18.22%"
$ws.Range("G21").Value = "4"
$ws.Range("D22").Value = "0.04097"
$ws.Range("E22").Value = "4.01%"
$ws.Range("G22").Value = "4"
$ws.Range("D23").Value = "0.001265"
$ws.Range("E23").Value = "0.29%"
$ws.Range("G23").Value = "4"
$ws.Range("D24").Value = "0.003891"
$ws.Range("E24").Value = "3.02%"
$ws.Range("G24").Value = "4"
$ws.Range("D25").Value = "0.0001232"
$ws.Range("E25").Value = "-3.96%"
$ws.Range("G25").Value = "4"
$ws.Range("G26").Value = "4"
$ws.Range("G27").Value = "4"
$ws.Range("G28").Value = "4"
$ws.Range("G29").Value = "4"
$ws.Range("G30").Value = "4"
$ws.Range("G31").Value = "4"
$ws.Range("G32").Value = "4"
$ws.Range("G33").Value = "4"
$ws.Range("G34").Value = "4"
$ws.Range("G35").Value = "4"
$ws.Range("G36").Value = "4"
$ws.Range("G37").Value = "4"
$ws.Range("D38").Value = "0.02415"
$ws.Range("E38").Value = "5.37%"
$ws.Range("G38").Value = "4"
$ws.Range("D39").Value = "0.05218"
$ws.Range("E39").Value = "4.86%"
$ws.Range("G39").Value = "4"
$ws.Range("D40").Value = "0.006304"
$ws.Range("E40").Value = "22.35%"
$ws.Range("G40").Value = "4"
$ws.Range("D41").Value = "0.007801"
$ws.Range("E41").Value = "1.55%"
$ws.Range("G41").Value = "4"
$ws.Range("E42").Value = "3.29%"
$ws.Range("G42").Value = "4"
$ws.Range("D43").Value = "0.007401"
$ws.Range("E43").Value = "0.56%"
$ws.Range("G43").Value = "4"
$ws.Range("D44").Value = "0.007425"
$ws.Range("E44").Value = "-6.34%"
$ws.Range("G44").Value = "4"
$ws.Range("D45").Value = "0.2955"
$ws.Range("E45").Value = "-6.43%"
$ws.Range("G45").Value = "4"
$ws.Range("D46").Value = "0.00006445"
$ws.Range("E46").Value = "1.36%"
$ws.Range("G46").Value = "4"
$ws.Range("E47").Value = "-0.11%"
$ws.Range("G47").Value = "4"
$ws.Range("D48").Value = "0.03342"
$ws.Range("E48").Value = "64.02%"
$ws.Range("G48").Value = "4"
$ws.Range("D49").Value = "0.004203"
$ws.Range("E49").Value = "0.04%"
$ws.Range("G49").Value = "4"
$ws.Range("D50").Value = "0.00002101"
$ws.Range("E50").Value = "-0.11%"
$ws.Range("G50").Value = "4"
$ws.Range("D51").Value = "0.0002001"
$ws.Range("E51").Value = "-0.11%"
$ws.Range("G51").Value = "4"

# Restore default (Normal) style on the touched numeric columns so no
# stray cell-format index is left behind (matches source formatting).
$colD.Style = "Normal"
$colE.Style = "Normal"
$colG.Style = "Normal"
